$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 206, shifting existing rows (206-235) down to (209-238)
$ws.Rows("206:208").Insert()

# New weekly data rows (same market/product metadata as neighboring rows, new date + figures)
$newRows = @(
    @{ Row = 206; K = "Hass"; L = "Primera"; M = 500; N = 23000; O = 26000; P = 24800; S = 2480; T = 10 },
    @{ Row = 207; K = "Hass"; L = "Segunda"; M = 150; N = 24000; O = 24000; P = 24000; S = 2400; T = 10 },
    @{ Row = 208; K = "Hass"; L = "Tercera"; M = 290; N = 21000; O = 22000; P = 21345; S = 2134; T = 10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($row, 4).Value = 45131
    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100106
    $ws.Cells.Item($row, 8).Value = "Oleaginosos"
    $ws.Cells.Item($row, 9).Value = 100106002
    $ws.Cells.Item($row, 10).Value = "Palta"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/bandeja 10 kilos"
    $ws.Cells.Item($row, 18).Value = "Perú"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}

Write-Host "Done"
